$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new "Basic Assumptions" sheet before "Historical Data"
# ---------------------------------------------------------------------------
$basic = $wb.Worksheets.Add($wb.Worksheets.Item("Historical Data"))
$basic.Name = "Basic Assumptions"
$hist = $wb.Worksheets.Item("Historical Data")

$basic.Range("A3").Value = "Currency"
$basic.Range("B3").Value = "USD"
$basic.Range("A4").Value = "Risk free rate"
$basic.Range("B4").Value = 0.044
$basic.Range("B4").NumberFormat = "0.00%"

# ---------------------------------------------------------------------------
# 2. Extend "Historical Data" with the 2023A / 2024A (L / M) columns
# ---------------------------------------------------------------------------

# Cost of sales (row 8)
$hist.Range("L8").Value = 150343
$hist.Range("M8").Value = 319099
$hist.Range("K8").Copy()
$hist.Range("L8:M8").PasteSpecial(-4122)

# R&D (row 9)
$hist.Range("L9").Value = 877387
$hist.Range("M9").Value = 804522
$hist.Range("K9").Copy()
$hist.Range("L9:M9").PasteSpecial(-4122)

# Selling, general and administrative (row 10)
$hist.Range("L10").Value = 481871
$hist.Range("M10").Value = 557872
$hist.Range("K10").Copy()
$hist.Range("L10:M10").PasteSpecial(-4122)

# Acquired in-process research and development (row 11) - n.a
$hist.Range("L11").Value = "n.a"
$hist.Range("M11").Value = "n.a"
$hist.Range("K11").Copy()
$hist.Range("L11:M11").PasteSpecial(-4122)

# Settlement and license charges (row 12) - n.a
$hist.Range("L12").Value = "n.a"
$hist.Range("M12").Value = "n.a"
$hist.Range("K12").Copy()
$hist.Range("L12:M12").PasteSpecial(-4122)

# Amortization of in-licensed rights (row 13)
$hist.Range("L13").Value = 1559
$hist.Range("M13").Value = 2405
$hist.Range("F13").Copy()
$hist.Range("L13:M13").PasteSpecial(-4122)

# Total Cost and Expenses (row 14) - sum formulas
$hist.Range("L14").Formula = "=SUM(L8:L13)"
$hist.Range("M14").Formula = "=SUM(M8:M13)"
$hist.Range("K14").Copy()
$hist.Range("L14:M14").PasteSpecial(-4122)

# Operating Loss (row 16) - difference formulas
$hist.Range("L16").Formula = "=L4-L14"
$hist.Range("M16").Formula = "=M4-M14"
$hist.Range("K16").Copy()
$hist.Range("L16:M16").PasteSpecial(-4122)

$hist.Range("K18").Select()

$basic.Activate()
$basic.Range("B5").Select()
